$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at Excel row 136 (shifts existing rows 136:153 down to 137:154)
$ws.Rows("136:136").Insert()

# Populate the newly inserted row 136 with the new data record
$ws.Cells.Item(136, 1).Value = 5
$ws.Cells.Item(136, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(136, 3).Value = 'Maule'
$ws.Cells.Item(136, 4).Value = 45212
$ws.Cells.Item(136, 5).Value = 7
$ws.Cells.Item(136, 6).Value = 100112013
$ws.Cells.Item(136, 7).Value = 'Alcachofa'
$ws.Cells.Item(136, 8).Value = 'Madrigal'
$ws.Cells.Item(136, 9).Value = 'Primera'
$ws.Cells.Item(136, 10).Value = 200
$ws.Cells.Item(136, 11).Value = 13000
$ws.Cells.Item(136, 12).Value = 13000
$ws.Cells.Item(136, 13).Value = 13000
$ws.Cells.Item(136, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(136, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(136, 16).Value = 325
$ws.Cells.Item(136, 17).Value = 40
$ws.Cells.Item(136, 18).Value = 'Hortaliza'
